{"js": "// Replace each arithmetic-problem cell text with its updated value.\n// Every \"old\" string is unique in the document, so a literal search\n// (matchCase, no wildcards) safely targets exactly one run each.\nconst body = context.document.body;\nconst pairs = [\n  [\"74+24=\", \"81-47=\"],\n  [\"77-20=\", \"9+20=\"],\n  [\"60-3=\", \"86-34=\"],\n  [\"9+41=\", \"97-57=\"],\n  [\"20+76=\", \"98-77=\"],\n  [\"33+21=\", \"31+53=\"],\n  [\"50+3=\", \"17-13=\"],\n  [\"39+60=\", \"86-32=\"],\n  [\"69+29=\", \"19+65=\"],\n  [\"36+6=\", \"14+77=\"],\n  [\"75+8=\", \"83+2=\"],\n  [\"61-11=\", \"95-26=\"],\n  [\"70+12=\", \"91-72=\"],\n  [\"63-63=\", \"44-1=\"],\n  [\"43+17=\", \"13+14=\"],\n  [\"82-39=\", \"17+71=\"],\n  [\"71-3=\", \"62-53=\"],\n  [\"86+10=\", \"93-66=\"],\n  [\"64+8=\", \"8+4=\"],\n  [\"71-39=\", \"92+6=\"],\n  [\"84-12=\", \"8+32=\"],\n  [\"66-64=\", \"92+7=\"],\n  [\"29-20=\", \"74+7=\"],\n  [\"42-30=\", \"12+15=\"],\n  [\"5+91=\", \"88-44=\"],\n  [\"60-35=\", \"70-23=\"],\n  [\"20+56=\", \"45+46=\"],\n  [\"11+71=\", \"56+27=\"],\n  [\"37-21=\", \"96-82=\"],\n  [\"7+57=\", \"74+2=\"],\n  [\"7+15=\", \"71-27=\"],\n  [\"92-84=\", \"96-21=\"],\n  [\"56+10=\", \"34-30=\"],\n  [\"62-37=\", \"26+40=\"],\n  [\"95-19=\", \"61-26=\"],\n  [\"33+6=\", \"94-47=\"],\n  [\"41-33=\", \"99-74=\"],\n  [\"68+17=\", \"30+52=\"],\n  [\"99-37=\", \"68+12=\"],\n  [\"48-25=\", \"93-7=\"],\n  [\"48+36=\", \"25+65=\"],\n  [\"45+5=\", \"65+31=\"],\n  [\"30+48=\", \"92-78=\"],\n  [\"17-1=\", \"30-2=\"],\n  [\"5+13=\", \"66-65=\"],\n  [\"5+70=\", \"67-10=\"],\n  [\"21-14=\", \"89-74=\"],\n  [\"57+7=\", \"49-1=\"],\n  [\"36+15=\", \"73+19=\"],\n  [\"0+65=\", \"73-46=\"],\n  [\"3+52=\", \"28-19=\"],\n  [\"81-53=\", \"67-60=\"],\n  [\"13+20=\", \"54-44=\"],\n  [\"18-3=\", \"53-5=\"],\n  [\"69+21=\", \"4+59=\"],\n  [\"4+75=\", \"63-21=\"],\n  [\"75-50=\", \"90-69=\"],\n  [\"89-38=\", \"69+1=\"],\n  [\"88-52=\", \"57-32=\"],\n  [\"47-31=\", \"95-42=\"],\n  [\"91-54=\", \"33-15=\"],\n  [\"8+50=\", \"16+74=\"],\n  [\"56-15=\", \"86-36=\"],\n  [\"55-30=\", \"74-70=\"],\n  [\"72-44=\", \"79-69=\"],\n  [\"37+24=\", \"90-71=\"],\n  [\"35+0=\", \"2+62=\"],\n  [\"35-34=\", \"92+6=\"],\n  [\"22+7=\", \"70-35=\"],\n  [\"4-2=\", \"60+39=\"],\n  [\"9+81=\", \"44-10=\"],\n  [\"59+36=\", \"21+45=\"],\n  [\"97-5=\", \"72-6=\"],\n  [\"81-3=\", \"12-1=\"],\n  [\"82+11=\", \"68-19=\"],\n  [\"90-50=\", \"48+10=\"],\n  [\"53+40=\", \"7+25=\"],\n  [\"51+46=\", \"80-57=\"],\n  [\"57-29=\", \"87-71=\"],\n  [\"64-17=\", \"4+19=\"],\n  [\"80-80=\", \"30+64=\"],\n  [\"13+63=\", \"89-62=\"],\n  [\"44+18=\", \"18+62=\"],\n  [\"89-54=\", \"6+61=\"],\n  [\"92-35=\", \"26-16=\"],\n  [\"14+14=\", \"47+12=\"],\n  [\"68-43=\", \"10+8=\"],\n  [\"14-1=\", \"85-61=\"],\n  [\"18+0=\", \"89-65=\"],\n  [\"46-34=\", \"40+46=\"],\n  [\"44+51=\", \"90-47=\"],\n  [\"73-70=\", \"71+10=\"],\n  [\"51+23=\", \"37-29=\"],\n  [\"88-1=\", \"36-1=\"],\n  [\"88-56=\", \"78-71=\"],\n  [\"49+45=\", \"2+31=\"],\n  [\"46-29=\", \"47-12=\"],\n  [\"55+8=\", \"57+38=\"],\n  [\"89-70=\", \"12+10=\"],\n  [\"11+21=\", \"83+9=\"]\n];\n\n// Phase 1: issue all searches, then resolve them with a single sync.\nconst resultsList = [];\nfor (const [oldText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  resultsList.push(results);\n}\nawait context.sync();\n\n// Phase 2: replace the text found by each search, then sync once more.\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const results = resultsList[i];\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace each arithmetic-problem cell text with its updated value.\n# Every \"old\" string is unique in the document, so Find/Replace (Replace = wdReplaceAll)\n# run one pair at a time safely targets exactly one run each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"74+24=\", \"81-47=\"),\n    @(\"77-20=\", \"9+20=\"),\n    @(\"60-3=\", \"86-34=\"),\n    @(\"9+41=\", \"97-57=\"),\n    @(\"20+76=\", \"98-77=\"),\n    @(\"33+21=\", \"31+53=\"),\n    @(\"50+3=\", \"17-13=\"),\n    @(\"39+60=\", \"86-32=\"),\n    @(\"69+29=\", \"19+65=\"),\n    @(\"36+6=\", \"14+77=\"),\n    @(\"75+8=\", \"83+2=\"),\n    @(\"61-11=\", \"95-26=\"),\n    @(\"70+12=\", \"91-72=\"),\n    @(\"63-63=\", \"44-1=\"),\n    @(\"43+17=\", \"13+14=\"),\n    @(\"82-39=\", \"17+71=\"),\n    @(\"71-3=\", \"62-53=\"),\n    @(\"86+10=\", \"93-66=\"),\n    @(\"64+8=\", \"8+4=\"),\n    @(\"71-39=\", \"92+6=\"),\n    @(\"84-12=\", \"8+32=\"),\n    @(\"66-64=\", \"92+7=\"),\n    @(\"29-20=\", \"74+7=\"),\n    @(\"42-30=\", \"12+15=\"),\n    @(\"5+91=\", \"88-44=\"),\n    @(\"60-35=\", \"70-23=\"),\n    @(\"20+56=\", \"45+46=\"),\n    @(\"11+71=\", \"56+27=\"),\n    @(\"37-21=\", \"96-82=\"),\n    @(\"7+57=\", \"74+2=\"),\n    @(\"7+15=\", \"71-27=\"),\n    @(\"92-84=\", \"96-21=\"),\n    @(\"56+10=\", \"34-30=\"),\n    @(\"62-37=\", \"26+40=\"),\n    @(\"95-19=\", \"61-26=\"),\n    @(\"33+6=\", \"94-47=\"),\n    @(\"41-33=\", \"99-74=\"),\n    @(\"68+17=\", \"30+52=\"),\n    @(\"99-37=\", \"68+12=\"),\n    @(\"48-25=\", \"93-7=\"),\n    @(\"48+36=\", \"25+65=\"),\n    @(\"45+5=\", \"65+31=\"),\n    @(\"30+48=\", \"92-78=\"),\n    @(\"17-1=\", \"30-2=\"),\n    @(\"5+13=\", \"66-65=\"),\n    @(\"5+70=\", \"67-10=\"),\n    @(\"21-14=\", \"89-74=\"),\n    @(\"57+7=\", \"49-1=\"),\n    @(\"36+15=\", \"73+19=\"),\n    @(\"0+65=\", \"73-46=\"),\n    @(\"3+52=\", \"28-19=\"),\n    @(\"81-53=\", \"67-60=\"),\n    @(\"13+20=\", \"54-44=\"),\n    @(\"18-3=\", \"53-5=\"),\n    @(\"69+21=\", \"4+59=\"),\n    @(\"4+75=\", \"63-21=\"),\n    @(\"75-50=\", \"90-69=\"),\n    @(\"89-38=\", \"69+1=\"),\n    @(\"88-52=\", \"57-32=\"),\n    @(\"47-31=\", \"95-42=\"),\n    @(\"91-54=\", \"33-15=\"),\n    @(\"8+50=\", \"16+74=\"),\n    @(\"56-15=\", \"86-36=\"),\n    @(\"55-30=\", \"74-70=\"),\n    @(\"72-44=\", \"79-69=\"),\n    @(\"37+24=\", \"90-71=\"),\n    @(\"35+0=\", \"2+62=\"),\n    @(\"35-34=\", \"92+6=\"),\n    @(\"22+7=\", \"70-35=\"),\n    @(\"4-2=\", \"60+39=\"),\n    @(\"9+81=\", \"44-10=\"),\n    @(\"59+36=\", \"21+45=\"),\n    @(\"97-5=\", \"72-6=\"),\n    @(\"81-3=\", \"12-1=\"),\n    @(\"82+11=\", \"68-19=\"),\n    @(\"90-50=\", \"48+10=\"),\n    @(\"53+40=\", \"7+25=\"),\n    @(\"51+46=\", \"80-57=\"),\n    @(\"57-29=\", \"87-71=\"),\n    @(\"64-17=\", \"4+19=\"),\n    @(\"80-80=\", \"30+64=\"),\n    @(\"13+63=\", \"89-62=\"),\n    @(\"44+18=\", \"18+62=\"),\n    @(\"89-54=\", \"6+61=\"),\n    @(\"92-35=\", \"26-16=\"),\n    @(\"14+14=\", \"47+12=\"),\n    @(\"68-43=\", \"10+8=\"),\n    @(\"14-1=\", \"85-61=\"),\n    @(\"18+0=\", \"89-65=\"),\n    @(\"46-34=\", \"40+46=\"),\n    @(\"44+51=\", \"90-47=\"),\n    @(\"73-70=\", \"71+10=\"),\n    @(\"51+23=\", \"37-29=\"),\n    @(\"88-1=\", \"36-1=\"),\n    @(\"88-56=\", \"78-71=\"),\n    @(\"49+45=\", \"2+31=\"),\n    @(\"46-29=\", \"47-12=\"),\n    @(\"55+8=\", \"57+38=\"),\n    @(\"89-70=\", \"12+10=\"),\n    @(\"11+21=\", \"83+9=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $result) {\n        Write-Output \"WARNING: no match for $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
